$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2998.05
$ws.Range("D3").Value = 402.89
$ws.Range("D4").Value = 36.28
$ws.Range("D5").Value = 149.21
$ws.Range("D6").Value = 50.74
$ws.Range("D7").Value = 137.51
$ws.Range("D8").Value = 858.8099999999999
$ws.Range("D9").Value = 918.87
$ws.Range("D10").Value = 3513.91
$ws.Range("D11").Value = 17.24
$ws.Range("D12").Value = 371.54
$ws.Range("D13").Value = 0.64
$ws.Range("D14").Value = 0.68
$ws.Range("D15").Value = 3.53
$ws.Range("D16").Value = 0.89
$ws.Range("D17").Value = 0.75
$ws.Range("D18").Value = 52.68
$ws.Range("D19").Value = 561.75
$ws.Range("D20").Value = 2.12
$ws.Range("D21").Value = 329.46
$ws.Range("D22").Value = 650.91
$ws.Range("D23").Value = 5.05
